# Update dq messages and test data (DQ_Report sheet + Statistik sheet)

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("DQ_Report")
$ws2 = $wb.Worksheets.Item("Statistik")

# ---- DQ_Report: rewrite data rows 2-21 (columns A-E) ----
# Columns: A = PatientIdentifikator, B = Aufnahmenummer, C = ICD_Primaerkode,
#          D = Orpha_Kode, E = dq_msg

$rows = @(
    @{ A = "P_20085651"; B = "F_101641"; C = "E75.0"; D = 846;    E = "Kodierung ist nicht eindeutig. Relation E75.0 - 846 ist im BfArM nicht vorhanden.  " },
    @{ A = "P_20085652"; B = "F_101642"; C = "E75.0"; D = 797;    E = "Kodierung ist nicht eindeutig. Relation E75.0 - 797 ist im BfArM nicht vorhanden.  " },
    @{ A = "P_20085653"; B = "F_101643"; C = "E75.0"; D = 309151; E = "Kodierung ist nicht eindeutig. Relation E75.0 - 309151 ist im BfArM nicht vorhanden.  " },
    @{ A = "P_20085654"; B = "F_101644"; C = "E75.0"; D = 309247; E = "Kodierung ist nicht eindeutig. Relation E75.0 - 309247 ist im BfArM nicht vorhanden.  " },
    @{ A = "P_20085751"; B = "F_101645"; C = "G70.0"; D = 586;    E = "Kodierung ist nicht eindeutig. Relation G70.0 - 586 ist im BfArM nicht vorhanden.  " },
    @{ A = "P_20085752"; B = "F_101646"; C = "I50.0"; D = 589;    E = "Kodierung ist nicht eindeutig. ICD10 Code I50.0 ist im BfArM Mapping nicht enthalten.  " },
    @{ A = "P_20085753"; B = "F_101647"; C = "E75.2"; D = 3;      E = "Kodierung ist nicht eindeutig. Relation E75.2 - 3 ist im BfArM nicht vorhanden.  " },
    @{ A = "P_20085757"; B = "F_101651"; C = $null;   D = 586;    E = "Fehlendes ICD10 Code.  " },
    @{ A = "P_20085758"; B = "F_101652"; C = $null;   D = 3;      E = "Kodierung ist nicht eindeutig. Orpha Code 3 ist im BfArM-Mapping nicht enthalten.  Fehlendes ICD10 Code.  " },
    @{ A = "P_20085761"; B = "F_101655"; C = $null;   D = $null;  E = "Fall ist nicht eindeutig. Fehlendes ICD10 Code.  " },
    @{ A = "P_20085762"; B = "F_101656"; C = "E66.89"; D = 320;   E = "Kodierung ist nicht eindeutig. Relation E66.89 - 320 ist im BfArM nicht vorhanden.  " },
    @{ A = "P_20085763"; B = "F_101657"; C = "G35.9"; D = 71529;  E = "Kodierung ist nicht eindeutig. Relation G35.9 - 71529 ist im BfArM nicht vorhanden.  " },
    @{ A = "P_20085765"; B = "F_101658"; C = "E75.2"; D = $null;  E = "ICD10 Kodierung E75.2 ist nicht eindeutig. ICD10-Orpha Relation ist gemäß Tracer-Diagnosenliste vom Typ 1-m.  Fehlendes Orpha_Kode.  " },
    @{ A = "P_20085766"; B = "F_101659"; C = "E75.0"; D = $null;  E = "ICD10 Kodierung E75.0 ist nicht eindeutig. ICD10-Orpha Relation ist gemäß Tracer-Diagnosenliste vom Typ 1-m.  Fehlendes Orpha_Kode.  " },
    @{ A = "P_20085767"; B = "F_101660"; C = "E74.0"; D = $null;  E = "ICD10 Kodierung E74.0 ist nicht eindeutig. ICD10-Orpha Relation ist gemäß Tracer-Diagnosenliste vom Typ 1-m.  Fehlendes Orpha_Kode.  " },
    @{ A = "P_20085768"; B = "F_101661"; C = "E75.2"; D = 342;    E = "Kodierung ist nicht eindeutig. Relation E75.2 - 342 ist im BfArM nicht vorhanden.  " },
    @{ A = "P_20085769"; B = "F_101662"; C = "E75.0"; D = 226;    E = "Kodierung ist nicht eindeutig. Relation E75.0 - 226 ist im BfArM nicht vorhanden.  " },
    @{ A = "P_20085772"; B = "F_101665"; C = "D45";   D = $null;  E = "Fehlendes Orpha_Kode.  " },
    @{ A = "P_20085773"; B = "F_101666"; C = "E84.0"; D = $null;  E = "Fehlendes Orpha_Kode.  " },
    @{ A = "P_20085774"; B = "F_101667"; C = "E84.1"; D = $null;  E = "Fehlendes Orpha_Kode.  " }
)

$r = 2
foreach ($row in $rows) {
    $ws1.Cells.Item($r, 1).Value = $row.A
    $ws1.Cells.Item($r, 2).Value = $row.B
    $ws1.Cells.Item($r, 3).Value = $row.C
    $ws1.Cells.Item($r, 4).Value = $row.D
    $ws1.Cells.Item($r, 5).Value = $row.E
    $r = $r + 1
}

# ---- Statistik: update aggregate metrics in row 2 ----
$ws2.Cells.Item(2, 2).Value = 44.29
$ws2.Cells.Item(2, 3).Value = 55.71
$ws2.Cells.Item(2, 4).Value = 68.42
$ws2.Cells.Item(2, 5).Value = 37.5
$ws2.Cells.Item(2, 6).Value = 19
$ws2.Cells.Item(2, 7).Value = 9
$ws2.Cells.Item(2, 8).Value = 27
$ws2.Cells.Item(2, 9).Value = 28
